# Applies the changes described in the commit:
#   - Inserts two new variable rows ("inccanc" and "dcens_canc") into the
#     "Variables" sheet (sheet1) right after "age_fup5" (i.e. before "bmi0"),
#     pushing every following row down by two positions. The two rows that
#     fall off the bottom of the previous range are simply the continuation
#     of the existing data (TEA_130302 / ART_SWEETENER_170201), which now
#     land on rows 78 and 79.
#   - Appends two new category rows for "inccanc" (0 = No, 1 = Yes) to the
#     "Categories" sheet (sheet2).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Variables
$ws2 = $wb.Worksheets.Item(2)   # Categories

# --- sheet1 ("Variables"): insert two rows before the current row 41 (bmi0) ---
$ws1.Range("A41:A42").EntireRow.Insert()

# New row 41: inccanc
$ws1.Range("A41").Value = 40
$ws1.Range("B41").Value = "inccanc"
$ws1.Range("C41").Value = "incident first occuring cancer"
$ws1.Range("D41").Value = "integer"

# New row 42: dcens_canc
$ws1.Range("A42").Value = 41
$ws1.Range("B42").Value = "dcens_canc"
$ws1.Range("C42").Value = "censored age for cancer at FUP5"
$ws1.Range("D42").Value = "date"

# The row-insert above also shifted the literal "index" values (column A)
# of every following row down along with their original row content, so
# column A no longer equals (row number - 1). Renumber column A for rows
# 43..79 to restore the simple index = row - 1 sequence.
for ($r = 43; $r -le 79; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 1
}

# --- sheet2 ("Categories"): append the two new category rows for inccanc ---
$ws2.Range("A73").Value = "inccanc"
$ws2.Range("B73").Value = 0
$ws2.Range("C73").Value = "No"

$ws2.Range("A74").Value = "inccanc"
$ws2.Range("B74").Value = 1
$ws2.Range("C74").Value = "Yes"
